$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.160.06'
$ws.Range("E2").Value = '  -6.29%  '

$ws.Range("D3").Value = '3.235.06'
$ws.Range("E3").Value = '  -7.77%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '511.15'
$ws.Range("E5").Value = '  -7.19%  '

$ws.Range("D6").Value = '167.75'
$ws.Range("E6").Value = '  -16.01%  '

$ws.Range("D7").Value = '0.588'
$ws.Range("E7").Value = '  -5.82%  '

$ws.Range("E8").Value = '  +0.16%  '

$ws.Range("D9").Value = '3.231.73'
$ws.Range("E9").Value = '  -7.57%  '

$ws.Range("D10").Value = '0.592'
$ws.Range("E10").Value = '  -9.21%  '

$ws.Range("D11").Value = '54.61'
$ws.Range("E11").Value = '  -11.64%  '

$ws.Range("D12").Value = '0.129'
$ws.Range("E12").Value = '  -9.73%  '

$ws.Range("D13").Value = '0.0000250'
$ws.Range("E13").Value = '  -6.73%  '

$ws.Range("D14").Value = '8.83'
$ws.Range("E14").Value = '  -9.86%  '

$ws.Range("D15").Value = '3.775.35'
$ws.Range("E15").Value = '  -7.38%  '

$ws.Range("D16").Value = '3.248.41'
$ws.Range("E16").Value = '  -7.70%  '

$ws.Range("E17").Value = '  -8.22%  '

$ws.Range("D18").Value = '63.110.78'
$ws.Range("E18").Value = '  -6.12%  '

$ws.Range("D19").Value = '16.96'
$ws.Range("E19").Value = '  -7.41%  '

$ws.Range("D20").Value = '10.77'
$ws.Range("E20").Value = '  -8.69%  '

$ws.Range("D21").Value = '0.936'
$ws.Range("E21").Value = '  -8.27%  '

$ws.Range("D22").Value = '363.52'
$ws.Range("E22").Value = '  -7.31%  '

$ws.Range("D23").Value = '3.65'
$ws.Range("E23").Value = '  -7.97%  '

$ws.Range("D24").Value = '78.47'
$ws.Range("E24").Value = '  -7.30%  '

$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D25").Value = '6.13'
$ws.Range("E25").Value = '  -0.47%  '

$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").Value = '10.65'
$ws.Range("E26").Value = '  -10.28%  '

$ws.Range("D27").Value = '3.83'
$ws.Range("E27").Value = '  -0.88%  '

$ws.Range("D28").Value = '2.60'
$ws.Range("E28").Value = '  -7.49%  '

$ws.Range("D29").Value = '11.04'
$ws.Range("E29").Value = '  -9.65%  '

$ws.Range("D30").Value = '8.06'
$ws.Range("E30").Value = '  -8.53%  '

$ws.Range("D31").Value = '28.04'
$ws.Range("E31").Value = '  -9.60%  '

$ws.Range("D32").Value = '627.91'
$ws.Range("E32").Value = '  -11.63%  '

$ws.Range("D33").Value = '6.43'
$ws.Range("E33").Value = '  -8.23%  '

$ws.Range("D34").Value = '10.99'
$ws.Range("E34").Value = '  -5.79%  '

$ws.Range("D35").Value = '58.15'
$ws.Range("E35").Value = '  -8.95%  '

$ws.Range("D36").Value = '0.102'
$ws.Range("E36").Value = '  -7.20%  '

$ws.Range("E37").Value = '  +0.01%  '

$ws.Range("D38").Value = '35.68'
$ws.Range("E38").Value = '  -6.70%  '

$ws.Range("D39").Value = '0.373'
$ws.Range("E39").Value = '  -5.21%  '

$ws.Range("E40").Value = '  +0.07%  '

$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '2.864.88'
$ws.Range("E41").Value = '  -6.57%  '

$ws.Range("B42").Value = 'PEPE'
$ws.Range("C42").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D42").Value = '0.0₃0665'
$ws.Range("E42").Value = '  -1.64%  '

$ws.Range("D43").Value = '0.119'
$ws.Range("E43").Value = '  -8.21%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = '2.61'
$ws.Range("E44").Value = '  -6.41%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = '2.37'
$ws.Range("E45").Value = '  -5.43%  '

$ws.Range("D46").Value = '2.60'
$ws.Range("E46").Value = '  -13.89%  '

$ws.Range("D47").Value = '0.0384'
$ws.Range("E47").Value = '  -5.51%  '

$ws.Range("D48").Value = '2.94'
$ws.Range("E48").Value = '  +1.48%  '

$ws.Range("D49").Value = '0.122'
$ws.Range("E49").Value = '  -5.39%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '133.73'
$ws.Range("E50").Value = '  -3.33%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '2.68'
$ws.Range("E51").Value = '  +2.17%  '
